$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data block (before former row 211),
# pushing all existing records down by two rows.
$ws.Rows.Item(211).Insert()
$ws.Rows.Item(211).Insert()

# New record #1 -> row 211
$ws.Cells.Item(211,1).Value2 = 10
$ws.Cells.Item(211,2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(211,3).Value2 = "La Araucanía"
$ws.Cells.Item(211,4).Value2 = 44505
$ws.Cells.Item(211,5).Value2 = 9
$ws.Cells.Item(211,6).Value2 = 100112037
$ws.Cells.Item(211,7).Value2 = "Cebollín"
$ws.Cells.Item(211,8).Value2 = "Sin especificar"
$ws.Cells.Item(211,9).Value2 = "Primera"
$ws.Cells.Item(211,10).Value2 = 65
$ws.Cells.Item(211,11).Value2 = 7000
$ws.Cells.Item(211,12).Value2 = 7000
$ws.Cells.Item(211,13).Value2 = 7000
$ws.Cells.Item(211,14).Value2 = "$/docena de paquetes"
$ws.Cells.Item(211,15).Value2 = "Provincia de Cautín"
$ws.Cells.Item(211,16).Value2 = 583
$ws.Cells.Item(211,17).Value2 = 12
$ws.Cells.Item(211,18).Value2 = "Hortaliza"

# New record #2 -> row 212
$ws.Cells.Item(212,1).Value2 = 10
$ws.Cells.Item(212,2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(212,3).Value2 = "La Araucanía"
$ws.Cells.Item(212,4).Value2 = 44505
$ws.Cells.Item(212,5).Value2 = 9
$ws.Cells.Item(212,6).Value2 = 100112037
$ws.Cells.Item(212,7).Value2 = "Cebollín"
$ws.Cells.Item(212,8).Value2 = "Sin especificar"
$ws.Cells.Item(212,9).Value2 = "Primera"
$ws.Cells.Item(212,10).Value2 = 55
$ws.Cells.Item(212,11).Value2 = 5000
$ws.Cells.Item(212,12).Value2 = 5000
$ws.Cells.Item(212,13).Value2 = 5000
$ws.Cells.Item(212,14).Value2 = "$/docena de paquetes"
$ws.Cells.Item(212,15).Value2 = "Región de O'Higgins"
$ws.Cells.Item(212,16).Value2 = 417
$ws.Cells.Item(212,17).Value2 = 12
$ws.Cells.Item(212,18).Value2 = "Hortaliza"

Write-Output "Inserted two new records at rows 211-212; dimension now A1:R243"
